# The edit swaps the two embedded themes in this deck: the slide master's
# theme (ppt/theme/theme1.xml, currently the "Integral" palette) is
# recoloured to the stock "Office Theme" palette that previously only lived
# on the notes-master theme (ppt/theme/theme2.xml). Only the colour scheme
# differs between the two theme parts (font scheme / format scheme are
# already identical), so re-pointing the 12 theme colour slots reproduces
# the target theme1.xml content.

function RGBVal($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$themeColors = $master.Theme.ThemeElements.ThemeColorScheme

# index -> (role, target "Office Theme" RGB)
$officeTheme = @(
    @(1,  0x00, 0x00, 0x00),  # dk1
    @(2,  0xFF, 0xFF, 0xFF),  # lt1
    @(3,  0x44, 0x54, 0x6A),  # dk2
    @(4,  0xE7, 0xE6, 0xE6),  # lt2
    @(5,  0x5B, 0x9B, 0xD5),  # accent1
    @(6,  0xED, 0x7D, 0x31),  # accent2
    @(7,  0xA5, 0xA5, 0xA5),  # accent3
    @(8,  0xFF, 0xC0, 0x00),  # accent4
    @(9,  0x44, 0x72, 0xC4),  # accent5
    @(10, 0x70, 0xAD, 0x47),  # accent6
    @(11, 0x05, 0x63, 0xC1),  # hlink
    @(12, 0x95, 0x4F, 0x72)   # folHlink
)

foreach ($entry in $officeTheme) {
    $idx = $entry[0]
    $r = $entry[1]
    $g = $entry[2]
    $b = $entry[3]
    $themeColors.Colors($idx).RGB = (RGBVal $r $g $b)
}
